$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update trade rows 2-18 to reflect the re-simulated backtest results
$data = @(
    @(1279, 82, 87, 0.03873328458999999, 0.04183080288549258, 3.961725899935013, 0.07997045250049073, 44594, 44599, 5),
    @(1168, 68, 90, 0.04239447212, 0.04429053450867629, 2.214600869973906, 0.04472428347048107, 44580, 44602, 22),
    @(1260, 124, 128, 0.03931961033, 0.04155490471642183, 2.816470926891507, 0.05684935246462386, 44636, 44640, 4),
    @(1168, 131, 138, 0.04240649412999999, 0.04728321038249495, 5.696004582914112, 0.1149992790619534, 44643, 44650, 7),
    @(1114, 110, 138, 0.04446562119999999, 0.04728321038249495, 3.138794349299384, 0.06336556437212115, 44622, 44650, 28),
    @(1111, 96, 138, 0.04458939484999999, 0.04728321038249495, 2.9928290566019, 0.0604138168180357, 44608, 44650, 42),
    @(1123, 89, 138, 0.04408704299999999, 0.04728321038249495, 3.589295970541839, 0.07249675108614029, 44601, 44650, 49),
    @(2440, 229, 240, 0.02030156128, 0.02105258700999914, 1.832502781197901, 0.0369934961967191, 44741, 44752, 11),
    @(2560, 243, 247, 0.01935061128, 0.02081059581652085, 3.737560413493384, 0.07544901374923674, 44755, 44759, 4),
    @(2328, 257, 261, 0.02127592467, 0.02333570310854898, 4.795164204942034, 0.09681264013184654, 44769, 44773, 4),
    @(2153, 264, 275, 0.02300891593, 0.02419069819893157, 2.544377225009672, 0.05136192737314982, 44776, 44787, 11),
    @(2115, 250, 275, 0.02342187848, 0.02419069819893157, 1.626053705540273, 0.03282485303593674, 44762, 44787, 25),
    @(2561, 341, 350, 0.01934676744, 0.02011052123410559, 1.955973466704432, 0.03947707525167821, 44853, 44862, 9),
    @(2465, 348, 360, 0.02009909902, 0.02084999676503673, 1.85096294151555, 0.03735977141510372, 44860, 44872, 12),
    @(2453, 306, 360, 0.02019379362, 0.02084999676503673, 1.609666314775108, 0.03249528827445425, 44818, 44872, 54),
    @(2416, 355, 365, 0.02050329281, 0.01706998, -8.294883748959997, -0.167451776737319, 44867, 44877, 10),
    @(2433, 327, 365, 0.02035815782, 0.01706998, -8.000136636059995, -0.1615164716313217, 44839, 44877, 38)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

